$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 4 (previously sub-2), shifting dimension to A1:AA3
$ws.Rows.Item(4).Delete()

# Update row 2 (sub-1) values
$ws.Range("C2").Value = "Control"
$ws.Range("D2").Value = "30y"
$ws.Range("F2").Value = "Human"
$ws.Range("G2").Value = "Not Defined"
$ws.Range("H2").Value = "Not Defined"
$ws.Range("I2").Value = "Prime Adult Stage"
$ws.Range("J2").Value = "Not Defined"
$ws.Range("K2").Value = "Not Defined"

# Update row 3 (sub-3 -> sub-2, and other values)
$ws.Range("A3").Value = "sub-2"
$ws.Range("C3").Value = "Control"
$ws.Range("D3").Value = "20y"
$ws.Range("F3").Value = "Human"
$ws.Range("G3").Value = "Not Defined"
$ws.Range("H3").Value = "Not Defined"
$ws.Range("I3").Value = "Prime Adult Stage"
$ws.Range("J3").Value = "Not Defined"
$ws.Range("K3").Value = "Not Defined"
